$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text semantics: a leading apostrophe tells Excel to treat the
# entry as literal text (matters for D-column values like "556.87" or
# "0.590" that would otherwise be auto-parsed as numbers and lose exact
# formatting), then reset Style back to Normal so no extra number-format /
# quote-prefix styling lingers on the cell (keeps cell style == original).

$ws.Range("D2").Value = "'62.435.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.45%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.433.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.91%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'556.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.26%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'139.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +7.82%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.80%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.432.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.97%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.43%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.58%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.28%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.65%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'26.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +12.62%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.864.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.77%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'62.283.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.30%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +8.04%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.431.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.37%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'11.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +7.86%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'347.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +11.65%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.75%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +3.97%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'65.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.53%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +2.20%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.25%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'1.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +14.94%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.20%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +15.42%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0₃0789"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.61%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +6.04%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'6.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +11.03%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'171.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.59%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +6.65%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +5.23%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'18.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.45%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'4.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +13.09%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'369.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +17.44%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.02%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.20%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +11.98%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'39.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.67%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'145.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.72%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'3.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.71%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'20.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +11.00%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Mantle"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.68%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Stellar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0954"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.15%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.33%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.97%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'17.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +7.41%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'BabyDogeCoin"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0₆0217"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.53%  "
$ws.Range("E51").Style = "Normal"
